$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: reorder the set string literal
$ws.Range("E2").Value = "{'list', 'List[any]', 'any'}"

# Row 6: move "Scalpel Accuracy:" label from C6 to E6, and its value from D6 to F6
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = "Scalpel Accuracy:"
$ws.Range("F6").Value = 100

# E7: fix wording from "over" to "vs"
$ws.Range("E7").Value = "Accuracy vs PyType"
